$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "dirty match" columns: E = uri params string, F = match flag text.
# Column F values ("true"/"false") must be stored as literal TEXT (shared string),
# not Excel's native Boolean type. A plain Value assignment of "true"/"false" gets
# auto-coerced to a Boolean, and a leading apostrophe (quote-prefix) would tag the
# cell style as text-quoted - neither matches the source data. Instead, write a
# text-literal formula (="true") and convert it to a static value via copy /
# paste-special-values, which keeps it a plain shared-string cell with no style change.
$ws.Range("E2").Value = "page=1&stream=false"
$ws.Range("F2").Formula = "=""true"""

$ws.Range("E3").Value = "page=1&stream=true"
$ws.Range("F3").Formula = "=""false"""

$ws.Range("E4").Value = "page=2&stream=false"
$ws.Range("F4").Formula = "=""true"""

$ws.Range("F2:F4").Copy()
$ws.Range("F2:F4").PasteSpecial(-4163)
